# edit.ps1
# 1) Bump the cached "datetimeFigureOut" field text (9/20/2017 -> 8/18/2018)
#    on the Slide Master and on every Custom Layout's Date placeholder.
# 2) On slide 21's code-sample textbox, merge the split "max<int> " / "(int x, int y) {"
#    runs into a single run, and the split "max<double> " / "(double " runs into a
#    single run (same visible text, fewer runs) to match the authored edit.

$p = $ppt.ActivePresentation

$oldDate = "9/20/2017"
$newDate = "8/18/2018"

function Update-DateFieldShape($shape) {
    if (-not $shape.HasTextFrame) { return }
    if ($shape.PlaceholderFormat.Type -ne 16) { return }   # 16 = ppPlaceholderDate
    $tr = $shape.TextFrame.TextRange
    if ($tr.Text -eq $oldDate) {
        $tr.Text = $newDate
    }
}

# Slide master's own Date placeholder.
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateFieldShape $master.Shapes.Item($i)
}

# Every custom (slide) layout's Date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        Update-DateFieldShape $layout.Shapes.Item($si)
    }
}

# --- Slide 21 code sample: merge split runs ---
$slide = $p.Slides.Item(21)
$codeShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.HasTextFrame -and ($sh.TextFrame.TextRange.Text.IndexOf("max<int>") -ge 0)) {
        $codeShape = $sh
    }
}

if ($codeShape -ne $null) {
    $tr = $codeShape.TextFrame.TextRange

    # "max<int> " + "(int x, int y) {" -> "max<int> (int x, int y) {"
    $oldSpan1 = "max<int> (int x, int y) {"
    $start1 = $tr.Text.IndexOf($oldSpan1) + 1
    $tr.Characters($start1, $oldSpan1.Length).Text = $oldSpan1

    # "max<double> " + "(double " -> "max<double> (double "
    $oldSpan2 = "max<double> (double "
    $start2 = $tr.Text.IndexOf($oldSpan2) + 1
    $tr.Characters($start2, $oldSpan2.Length).Text = $oldSpan2
}
